# Adds the three new "QUOTA_EXP_GRPx" group-semaphore indicator rows to the
# "Variable" sheet and the matching "..._REAL" data-type rows to the
# "r Variable_DataType" sheet (aggiunti indicatori per semaforo di gruppo).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Variable": append rows 164-166
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Variable")

$row = 164
$ids        = @("INDICATOR_QUOTA_EXP_GRP1","INDICATOR_QUOTA_EXP_GRP2","INDICATOR_QUOTA_EXP_GRP3")
$columnName = @("QUOTA_EXP_GRP1","QUOTA_EXP_GRP2","QUOTA_EXP_GRP3")
$columnLbl  = @("500 - QUOTA_EXP_GRP1","501 - QUOTA_EXP_GRP2","503 - QUOTA_EXP_GRP3")
$colOrder   = @(500,501,502)

# Column B/C ("Id"/"Name") first, for every new row, so the shared strings
# are interned in the same order the original author typed them.
for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws1.Cells.Item($r,2).Value = $ids[$i]
    $ws1.Cells.Item($r,3).Value = $ids[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws1.Cells.Item($r,5).Value = $columnName[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws1.Cells.Item($r,6).Value = $columnLbl[$i]
}

# Remaining columns (A, D, G..N)
for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws1.Cells.Item($r,1).Value = "CREATE/MODIFY"
    $ws1.Cells.Item($r,7).Value = $colOrder[$i]
    $ws1.Cells.Item($r,8).Value = "false"
    $ws1.Cells.Item($r,9).Value = $true
    $ws1.Cells.Item($r,10).Value = "CONTINUOUS"
    $ws1.Cells.Item($r,11).Value = "false"
    $ws1.Cells.Item($r,12).Value = $false
    $ws1.Cells.Item($r,13).Value = $false
    $ws1.Cells.Item($r,14).Value = $false
}

# Highlight the new rows (A:N) in yellow, matching the author's manual markup.
$ws1.Range("A164:N166").Interior.Color = 65535
# Column D has no content on the new rows but keeps the usual left/top aligned,
# bordered look used throughout the sheet.
$ws1.Range("D164:D166").HorizontalAlignment = -4131
$ws1.Range("D164:D166").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# Sheet "r Variable_DataType": append rows 164-166
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("r Variable_DataType")

$realIds = @("INDICATOR_QUOTA_EXP_GRP1_REAL","INDICATOR_QUOTA_EXP_GRP2_REAL","INDICATOR_QUOTA_EXP_GRP3_REAL")

for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws2.Cells.Item($r,2).Value = $realIds[$i]
    $ws2.Cells.Item($r,3).Value = $realIds[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $r = $row + $i
    $ws2.Cells.Item($r,1).Value = "CREATE/MODIFY"
    $ws2.Cells.Item($r,5).Value = $ids[$i]
    $ws2.Cells.Item($r,6).Value = "REAL"
}

$ws2.Range("A164:E166").Interior.Color = 65535

# ---------------------------------------------------------------------------
# Leave the cursor roughly where the author left it.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("W144").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B171").Select() | Out-Null

$ws1.Activate()
